# Apply cryptos list update (prices and 1h volume/% changes) per commit
# "Updated cryptos list on Wed Apr 17 15:04:08 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.263.00"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "2.981.78"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'528.88"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "'131.43"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "2.978.95"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").Value = "'0.486"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "'0.149"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'6.15"
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").Value = "'0.441"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "'0.0000217"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("D14").Value = "'33.39"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "3.474.69"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.110"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "61.347.51"
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("D18").Value = "2.984.78"
$ws.Range("E18").Value = "  -2.39%  "
$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'457.96"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "'13.16"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'0.672"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "'6.84"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").Value = "'77.83"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'11.77"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'2.64"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "'7.63"
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'25.41"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +5.00%  "
$ws.Range("D32").Value = "'1.83"
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").Value = "'55.44"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "'5.38"
$ws.Range("E34").Value = "  +5.12%  "
$ws.Range("D35").Value = "'2.23"
$ws.Range("E35").Value = "  -5.20%  "
$ws.Range("D36").Value = "'5.77"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "'457.95"
$ws.Range("E37").Value = "  -2.20%  "
$ws.Range("D38").Value = "3.153.13"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").Value = "'0.0383"
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("D40").Value = "'0.0778"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'0.119"
$ws.Range("E41").Value = "  +7.13%  "
$ws.Range("D42").Value = "'7.98"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").Value = "'2.38"
$ws.Range("E43").Value = "  -6.87%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'0.244"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "'25.38"
$ws.Range("E46").Value = "  +6.33%  "
$ws.Range("D47").Value = "'120.44"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("D49").Value = "0.0₃0505"
$ws.Range("E49").Value = "  -0.67%  "
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("E51").Value = "  +6.24%  "
